$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.401.50'
$ws.Range("E2").Value = '  -0.06%  '
$ws.Range("D3").Value = '1.843.94'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9990'
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '239.13'
$ws.Range("E5").Value = '  -0.80%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6312'
$ws.Range("E6").Value = '  -0.25%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.000'
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07527'
$ws.Range("E8").Value = '  -0.41%  '
$ws.Range("E9").Value = '  -1.51%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '24.43'
$ws.Range("E10").Value = '  -1.10%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07716'
$ws.Range("E11").Value = '  +0.09%  '
$ws.Range("D12").Value = '1.843.09'
$ws.Range("E12").Value = '  -7.15%  '
$ws.Range("E13").Value = '  +0.25%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6796'
$ws.Range("E14").Value = '  -0.87%  '
$ws.Range("E15").Value = '  +4.41%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '83.15'
$ws.Range("E16").Value = '  +0.06%  '
$ws.Range("D17").Value = '2.099.79'
$ws.Range("E17").Value = '  -7.27%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.167'
$ws.Range("E18").Value = '  -0.19%  '
$ws.Range("D19").Value = '29.426.59'
$ws.Range("E19").Value = '  -0.05%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '228.28'
$ws.Range("E20").Value = '  -1.55%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.40'
$ws.Range("E21").Value = '  -0.77%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.0000'
$ws.Range("E23").Value = '  -1.81%  '
$ws.Range("E24").Value = '  +0.06%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '157.13'
$ws.Range("E25").Value = '  +1.33%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1392'
$ws.Range("E26").Value = '  +0.35%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.373'
$ws.Range("E27").Value = '  -0.46%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '17.59'
$ws.Range("E28").Value = '  -0.47%  '
$ws.Range("E29").Value = '  -0.65%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.278'
$ws.Range("E30").Value = '  +1.62%  '
$ws.Range("E31").Value = '  -2.87%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.101'
$ws.Range("E32").Value = '  -0.80%  '
$ws.Range("E33").Value = '  +0.15%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.842'
$ws.Range("E34").Value = '  -1.76%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.158'
$ws.Range("E35").Value = '  +0.10%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7091'
$ws.Range("E36").Value = '  -1.10%  '
$ws.Range("E37").Value = '  -0.09%  '
$ws.Range("D38").Value = '1.247.43'
$ws.Range("E38").Value = '  -0.28%  '
$ws.Range("E39").Value = '  +0.05%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.765'
$ws.Range("E40").Value = '  -1.05%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.323'
$ws.Range("E41").Value = '  +3.78%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.9025'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.0000'
$ws.Range("E43").Value = '  +0.08%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '101.85'
$ws.Range("E44").Value = '  +0.12%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '65.82'
$ws.Range("E45").Value = '  -1.74%  '
$ws.Range("B46").Value = 'BabyDogeCoin'
$ws.Range("C46").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00000000119'
$ws.Range("E46").Value = '  +0.56%  '
$ws.Range("B47").Value = 'Aptos'
$ws.Range("C47").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '7.096'
$ws.Range("E47").Value = '  -1.53%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.3995'
$ws.Range("E48").Value = '  -0.54%  '
$ws.Range("E49").Value = '  -1.07%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '8.901'
$ws.Range("E50").Value = '  -3.22%  '
$ws.Range("E51").Value = '  -0.48%  '
